# Apply crypto price/volume updates scraped on Fri Mar  3 14:54:52 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "22.411.02"
$ws.Range("E2").Value = "  -3.99%  "
# Row 3
$ws.Range("D3").Value = "1.573.53"
$ws.Range("E3").Value = "  -3.43%  "
# Row 4
$ws.Range("E4").Value = "  -0.10%  "
# Row 5
$ws.Range("D5").Value = "'1.001"
$ws.Range("E5").Value = "  -0.06%  "
# Row 6
$ws.Range("D6").Value = "'289.68"
$ws.Range("E6").Value = "  -2.61%  "
# Row 7
$ws.Range("D7").Value = "'0.3677"
$ws.Range("E7").Value = "  -2.22%  "
# Row 8
$ws.Range("D8").Value = "'49.40"
$ws.Range("E8").Value = "  -1.05%  "
# Row 9
$ws.Range("D9").Value = "'0.3397"
$ws.Range("E9").Value = "  -2.96%  "
# Row 10
$ws.Range("D10").Value = "'1.170"
$ws.Range("E10").Value = "  -2.56%  "
# Row 11
$ws.Range("D11").Value = "'0.07624"
$ws.Range("E11").Value = "  -4.85%  "
# Row 12
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  -0.13%  "
# Row 13
$ws.Range("D13").Value = "'21.42"
$ws.Range("E13").Value = "  -1.85%  "
# Row 14
$ws.Range("D14").Value = "'6.070"
$ws.Range("E14").Value = "  -3.33%  "
# Row 15
$ws.Range("D15").Value = "'6.940"
$ws.Range("E15").Value = "  -3.84%  "
# Row 16
$ws.Range("D16").Value = "'0.00001139"
$ws.Range("E16").Value = "  -4.12%  "
# Row 17
$ws.Range("D17").Value = "1.566.30"
$ws.Range("E17").Value = "  -4.17%  "
# Row 18
$ws.Range("D18").Value = "'89.38"
$ws.Range("E18").Value = "  -5.89%  "
# Row 19
$ws.Range("D19").Value = "'0.06762"
$ws.Range("E19").Value = "  -2.50%  "
# Row 20
$ws.Range("E20").Value = "  -0.14%  "
# Row 21
$ws.Range("D21").Value = "'6.259"
$ws.Range("E21").Value = "  -5.72%  "
# Row 22
$ws.Range("D22").Value = "'16.60"
$ws.Range("E22").Value = "  -4.01%  "
# Row 23
$ws.Range("D23").Value = "'0.5294"
$ws.Range("E23").Value = "  -7.05%  "
# Row 24
$ws.Range("D24").Value = "'12.01"
$ws.Range("E24").Value = "  -1.82%  "
# Row 25
$ws.Range("D25").Value = "22.414.75"
$ws.Range("E25").Value = "  -4.04%  "
# Row 26
$ws.Range("D26").Value = "'2.381"
$ws.Range("E26").Value = "  -2.33%  "
# Row 27
$ws.Range("D27").Value = "'2.962"
$ws.Range("E27").Value = "  +1.08%  "
# Row 28
$ws.Range("D28").Value = "'20.02"
$ws.Range("E28").Value = "  -3.41%  "
# Row 29
$ws.Range("D29").Value = "'146.18"
$ws.Range("E29").Value = "  -3.35%  "
# Row 30
$ws.Range("D30").Value = "'4.978"
$ws.Range("E30").Value = "  -3.63%  "
# Row 31
$ws.Range("D31").Value = "'125.96"
$ws.Range("E31").Value = "  -4.35%  "
# Row 32
$ws.Range("D32").Value = "1.741.31"
$ws.Range("E32").Value = "  -4.06%  "
# Row 33
$ws.Range("D33").Value = "'1.047"
$ws.Range("E33").Value = "  +8.80%  "
# Row 34
$ws.Range("D34").Value = "'6.298"
$ws.Range("E34").Value = "  -6.77%  "
# Row 35
$ws.Range("E35").Value = "  -5.06%  "
# Row 36
$ws.Range("E36").Value = "  -7.17%  "
# Row 37
$ws.Range("D37").Value = "'0.08452"
$ws.Range("E37").Value = "  -2.75%  "
# Row 38
$ws.Range("D38").Value = "'0.02543"
$ws.Range("E38").Value = "  -5.07%  "
# Row 39
$ws.Range("D39").Value = "'0.2333"
$ws.Range("E39").Value = "  -3.33%  "
# Row 40
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "'5.575"
$ws.Range("E40").Value = "  -4.16%  "
# Row 41
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "'0.06570"
$ws.Range("E41").Value = "  -2.53%  "
# Row 42
$ws.Range("E42").Value = "  -7.54%  "
# Row 43
$ws.Range("D43").Value = "'1.249"
$ws.Range("E43").Value = "  -3.31%  "
# Row 44
$ws.Range("D44").Value = "'0.6389"
$ws.Range("E44").Value = "  -5.89%  "
# Row 45
$ws.Range("D45").Value = "'14.34"
$ws.Range("E45").Value = "  -6.22%  "
# Row 46
$ws.Range("E46").Value = "  -0.08%  "
# Row 47
$ws.Range("D47").Value = "'0.6014"
$ws.Range("E47").Value = "  -4.36%  "
# Row 48
$ws.Range("D48").Value = "'3.744"
$ws.Range("E48").Value = "  -3.56%  "
# Row 49
$ws.Range("D49").Value = "'2.126"
$ws.Range("E49").Value = "  -4.50%  "
# Row 50
$ws.Range("D50").Value = "'1.259"
$ws.Range("E50").Value = "  +4.02%  "
# Row 51
$ws.Range("D51").Value = "'123.60"
$ws.Range("E51").Value = "  -2.21%  "
